$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: 45181 | 4 botellones | -212
$ws.Range("B32").Value = 45181
$ws.Range("C32").Value = "4 botellones"
$ws.Range("D32").Value = -212

# Row 33: 45184 | 4 botellones | -212
$ws.Range("B33").Value = 45184
$ws.Range("C33").Value = "4 botellones"
$ws.Range("D33").Value = -212

# Row 34: 45189 | a comprobar | 500
$ws.Range("B34").Value = 45189
$ws.Range("C34").Value = "a comprobar"
$ws.Range("D34").Value = 500

# Update the active selection to match the authored state
$ws.Range("D35").Select()
